# Case with 380 kV: update loading_percent values (Sheet1, rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 30.59289261161933
$ws.Range("C2").Value = 8.349247353293682
$ws.Range("D2").Value = 13.29410330362165
$ws.Range("E2").Value = 12.75344722785181
$ws.Range("G2").Value = 3.890146777832082
$ws.Range("J2").Value = 7.845001912645786
$ws.Range("K2").Value = 25.49035728962296
$ws.Range("L2").Value = 14.60392309422107
$ws.Range("N2").Value = 31.65524084359587

# Row 3
$ws.Range("B3").Value = 30.51139263075643
$ws.Range("C3").Value = 8.252932480579201
$ws.Range("D3").Value = 13.3002580987778
$ws.Range("E3").Value = 12.77764170490774
$ws.Range("G3").Value = 3.894418003580121
$ws.Range("J3").Value = 7.851691727538383
$ws.Range("K3").Value = 25.43914056227185
$ws.Range("L3").Value = 14.62300235898947
$ws.Range("N3").Value = 31.58472876753818

# Row 4
$ws.Range("B4").Value = 30.46828476583568
$ws.Range("C4").Value = 8.195789625228258
$ws.Range("D4").Value = 13.30641649061195
$ws.Range("E4").Value = 12.79403516926478
$ws.Range("G4").Value = 3.897175179570771
$ws.Range("J4").Value = 7.856029908973998
$ws.Range("K4").Value = 25.41331238699426
$ws.Range("L4").Value = 14.636989355413
$ws.Range("N4").Value = 31.5423533303027

# Row 5
$ws.Range("B5").Value = 30.45247087250158
$ws.Range("C5").Value = 8.173030509015131
$ws.Range("D5").Value = 13.30952421520336
$ws.Range("E5").Value = 12.8011026342866
$ws.Range("G5").Value = 3.898332738827066
$ws.Range("J5").Value = 7.857855935070819
$ws.Range("K5").Value = 25.40420528535899
$ws.Range("L5").Value = 14.64326024447635
$ws.Range("N5").Value = 31.52531889579169

# Row 6
$ws.Range("B6").Value = 30.44995109994986
$ws.Range("C6").Value = 8.169283927378222
$ws.Range("D6").Value = 13.31007636729598
$ws.Range("E6").Value = 12.8022995625762
$ws.Range("G6").Value = 3.898527007386584
$ws.Range("J6").Value = 7.85816266505507
$ws.Range("K6").Value = 25.40277883248846
$ws.Range("L6").Value = 14.644336001844
$ws.Range("N6").Value = 31.52250454458248

# Row 7
$ws.Range("B7").Value = 30.46806438562832
$ws.Range("C7").Value = 8.195480521086571
$ws.Range("D7").Value = 13.30645598103429
$ws.Range("E7").Value = 12.79412891629992
$ws.Range("G7").Value = 3.897190653028492
$ws.Range("J7").Value = 7.856054299567106
$ws.Range("K7").Value = 25.41318381788575
$ws.Range("L7").Value = 14.63707161503236
$ws.Range("N7").Value = 31.54212264755137

# Row 8
$ws.Range("B8").Value = 30.56335856562711
$ws.Range("C8").Value = 8.315639572866225
$ws.Range("D8").Value = 13.29573159959577
$ws.Range("E8").Value = 12.76147044072993
$ws.Range("G8").Value = 3.891591636161208
$ws.Range("J8").Value = 7.847260814778553
$ws.Range("K8").Value = 25.47153495026632
$ws.Range("L8").Value = 14.61002985990626
$ws.Range("N8").Value = 31.63073916868115

# Row 9
$ws.Range("B9").Value = 30.80475322364673
$ws.Range("C9").Value = 8.565927983703615
$ws.Range("D9").Value = 13.29358148648431
$ws.Range("E9").Value = 12.70961884643243
$ws.Range("G9").Value = 3.881674027790798
$ws.Range("J9").Value = 7.831837804509567
$ws.Range("K9").Value = 25.6302700131874
$ws.Range("L9").Value = 14.57504204964618
$ws.Range("N9").Value = 31.81168884080631

# Row 10
$ws.Range("B10").Value = 31.01457901786749
$ws.Range("C10").Value = 8.757123013593741
$ws.Range("D10").Value = 13.30350488935062
$ws.Range("E10").Value = 12.67894061959297
$ws.Range("G10").Value = 3.875026402050967
$ws.Range("J10").Value = 7.821604418772258
$ws.Range("K10").Value = 25.77346027375544
$ws.Range("L10").Value = 14.56034810048766
$ws.Range("N10").Value = 31.94884281180881

# Row 11
$ws.Range("B11").Value = 31.11688720729466
$ws.Range("C11").Value = 8.845333323512921
$ws.Range("D11").Value = 13.31051146869408
$ws.Range("E11").Value = 12.66659150745999
$ws.Range("G11").Value = 3.872139101084102
$ws.Range("J11").Value = 7.817184786613645
$ws.Range("K11").Value = 25.84425041955632
$ws.Range("L11").Value = 14.55605560164532
$ws.Range("N11").Value = 32.01211428782857

# Row 12
$ws.Range("B12").Value = 31.15659533672284
$ws.Range("C12").Value = 8.878883792687716
$ws.Range("D12").Value = 13.31352217568674
$ws.Range("E12").Value = 12.66214595004076
$ws.Range("G12").Value = 3.871065275965364
$ws.Range("J12").Value = 7.815544864602971
$ws.Range("K12").Value = 25.87185747019494
$ws.Range("L12").Value = 14.55477396796821
$ws.Range("N12").Value = 32.03619630795483

# Row 13
$ws.Range("B13").Value = 31.14800083899879
$ws.Range("C13").Value = 8.871652037005404
$ws.Range("D13").Value = 13.31285788483404
$ws.Range("E13").Value = 12.66309312016163
$ws.Range("G13").Value = 3.871295676737103
$ws.Range("J13").Value = 7.81589655533458
$ws.Range("K13").Value = 25.86587641148672
$ws.Range("L13").Value = 14.55503470097773
$ws.Range("N13").Value = 32.03100443151781

# Row 14
$ws.Range("B14").Value = 31.12013477323089
$ws.Range("C14").Value = 8.848090765457108
$ws.Range("D14").Value = 13.31075200332467
$ws.Range("E14").Value = 12.6662211449475
$ws.Range("G14").Value = 3.872050366114197
$ws.Range("J14").Value = 7.817049194937733
$ws.Range("K14").Value = 25.84650570301362
$ws.Range("L14").Value = 14.5559432707184
$ws.Range("N14").Value = 32.01409309870613

# Row 15
$ws.Range("B15").Value = 31.10319122104988
$ws.Range("C15").Value = 8.833677068574865
$ws.Range("D15").Value = 13.30950861142215
$ws.Range("E15").Value = 12.66816719794401
$ws.Range("G15").Value = 3.872515175095333
$ws.Range("J15").Value = 7.81775960296121
$ws.Range("K15").Value = 25.83474442711702
$ws.Range("L15").Value = 14.55654456951598
$ws.Range("N15").Value = 32.00375024937125

# Row 16
$ws.Range("B16").Value = 31.00802926395442
$ws.Range("C16").Value = 8.75138050515088
$ws.Range("D16").Value = 13.30309707119186
$ws.Range("E16").Value = 12.67977997119921
$ws.Range("G16").Value = 3.875217834378898
$ws.Range("J16").Value = 7.821897977290551
$ws.Range("K16").Value = 25.76894667219826
$ws.Range("L16").Value = 14.56067674914785
$ws.Range("N16").Value = 31.94472523248217

# Row 17
$ws.Range("B17").Value = 30.9513931988343
$ws.Range("C17").Value = 8.701188795616229
$ws.Range("D17").Value = 13.29980156262224
$ws.Range("E17").Value = 12.6873153456404
$ws.Range("G17").Value = 3.876910757705859
$ws.Range("J17").Value = 7.82449694656784
$ws.Range("K17").Value = 25.7300210696984
$ws.Range("L17").Value = 14.56382427806217
$ws.Range("N17").Value = 31.90873776614875

# Row 18
$ws.Range("B18").Value = 30.91946437017026
$ws.Range("C18").Value = 8.672437888963215
$ws.Range("D18").Value = 13.29814073727475
$ws.Range("E18").Value = 12.6918007265304
$ws.Range("G18").Value = 3.877897360516757
$ws.Range("J18").Value = 7.826013990163819
$ws.Range("K18").Value = 25.7081649201239
$ws.Range("L18").Value = 14.56585979871393
$ws.Range("N18").Value = 31.88812175376464

# Row 19
$ws.Range("B19").Value = 30.90876547451438
$ws.Range("C19").Value = 8.662724527849578
$ws.Range("D19").Value = 13.29761873887506
$ws.Range("E19").Value = 12.69334538087875
$ws.Range("G19").Value = 3.878233623052072
$ws.Range("J19").Value = 7.826531450899528
$ws.Range("K19").Value = 25.70085667650566
$ws.Range("L19").Value = 14.56658766135376
$ws.Range("N19").Value = 31.88115590451099

# Row 20
$ws.Range("B20").Value = 30.95735540860686
$ws.Range("C20").Value = 8.706519793194895
$ws.Range("D20").Value = 13.30012809435566
$ws.Range("E20").Value = 12.68649754184185
$ws.Range("G20").Value = 3.876729211234413
$ws.Range("J20").Value = 7.824217987077074
$ws.Range("K20").Value = 25.7341097072486
$ws.Range("L20").Value = 14.56346591697605
$ws.Range("N20").Value = 31.91256012875313

# Row 21
$ws.Range("B21").Value = 31.12829365449772
$ws.Range("C21").Value = 8.855007530187146
$ws.Range("D21").Value = 13.31136085874684
$ws.Range("E21").Value = 12.66529610681719
$ws.Range("G21").Value = 3.871828166458207
$ws.Range("J21").Value = 7.816709723655905
$ws.Range("K21").Value = 25.85217373272891
$ws.Range("L21").Value = 14.55566707131494
$ws.Range("N21").Value = 32.01905707538183

# Row 22
$ws.Range("B22").Value = 31.24563339639838
$ws.Range("C22").Value = 8.952896558491615
$ws.Range("D22").Value = 13.32078512353705
$ws.Range("E22").Value = 12.65278480550535
$ws.Range("G22").Value = 3.868738844331548
$ws.Range("J22").Value = 7.811998967225723
$ws.Range("K22").Value = 25.93399395098204
$ws.Range("L22").Value = 14.55257412995878
$ws.Range("N22").Value = 32.08937096475581

# Row 23
$ws.Range("B23").Value = 31.182499525419
$ws.Range("C23").Value = 8.900584132104557
$ws.Range("D23").Value = 13.31556498677283
$ws.Range("E23").Value = 12.6593393314032
$ws.Range("G23").Value = 3.870377304247848
$ws.Range("J23").Value = 7.814495281830881
$ws.Range("K23").Value = 25.88990309603548
$ws.Range("L23").Value = 14.55404158088514
$ws.Range("N23").Value = 32.05177932885348

# Row 24
$ws.Range("B24").Value = 30.95465792354459
$ws.Range("C24").Value = 8.704109319904644
$ws.Range("D24").Value = 13.29997974103578
$ws.Range("E24").Value = 12.6868667937571
$ws.Range("G24").Value = 3.876811246899533
$ws.Range("J24").Value = 7.824344033440784
$ws.Range("K24").Value = 25.73225960519711
$ws.Range("L24").Value = 14.56362722818798
$ws.Range("N24").Value = 31.91083180888112

# Row 25
$ws.Range("B25").Value = 30.73368467253154
$ws.Range("C25").Value = 8.496814819678068
$ws.Range("D25").Value = 13.29214174541767
$ws.Range("E25").Value = 12.72234240998645
$ws.Range("G25").Value = 3.884244195915843
$ws.Range("J25").Value = 7.83581645623048
$ws.Range("K25").Value = 25.5826253732754
$ws.Range("L25").Value = 14.58257401279878
$ws.Range("N25").Value = 31.76198531851732

